$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 3.1
$ws.Range("L2").Value = 3.1
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("W2").Value = 7.5
$ws.Range("AE2").Value = 21
$ws.Range("AP2").Value = 41
$ws.Range("AT2").Value = 2.2
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("AH9").Value = 10
$ws.Range("AN9").Value = 5
$ws.Range("AX9").Value = 13
$ws.Range("G10").Value = 2.5
$ws.Range("I10").Value = 3.1
$ws.Range("L10").Value = 4
$ws.Range("W10").Value = 6
$ws.Range("X10").Value = 10
$ws.Range("Z10").Value = 23
$ws.Range("AA10").Value = 23
$ws.Range("AR10").Value = 81
$ws.Range("G14").Value = 5.5
$ws.Range("I14").Value = 1.53
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 9
$ws.Range("O14").Value = 1.33
$ws.Range("P14").Value = 3.25
$ws.Range("R14").Value = 1.75
$ws.Range("AC14").Value = 9
$ws.Range("G15").Value = 3.6
$ws.Range("H15").Value = 2.9
$ws.Range("J15").Value = 4.5
$ws.Range("M15").Value = 1.1
$ws.Range("N15").Value = 7
$ws.Range("O15").Value = 1.5
$ws.Range("P15").Value = 2.5
$ws.Range("Q15").Value = 2.5
$ws.Range("R15").Value = 1.5
$ws.Range("S15").Value = 1.57
$ws.Range("T15").Value = 2.25
$ws.Range("U15").Value = 2.1
$ws.Range("V15").Value = 1.67
$ws.Range("AC15").Value = 6.5
$ws.Range("AE15").Value = 19
$ws.Range("AH15").Value = 9
$ws.Range("AR15").Value = 126
$ws.Range("AS15").Value = 401
$ws.Range("AT15").Value = 2.25
$ws.Range("BD15").Value = 151
$ws.Range("G16").Value = 2.88
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 2.3
$ws.Range("J16").Value = 3.5
$ws.Range("L16").Value = 3
$ws.Range("S16").Value = 1.4
$ws.Range("T16").Value = 2.75
$ws.Range("U16").Value = 1.73
$ws.Range("V16").Value = 2
$ws.Range("Z16").Value = 29
$ws.Range("AJ16").Value = 23
$ws.Range("AK16").Value = 19
$ws.Range("AM16").Value = 201
$ws.Range("AT16").Value = 2.75
$ws.Range("AU16").Value = 8
$ws.Range("AW16").Value = 4.5
$ws.Range("AY16").Value = 23
$ws.Range("BA16").Value = 67
$ws.Range("G18").Value = 7.5
$ws.Range("H18").Value = 5.1
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = 2.72
$ws.Range("S18").Value = 1.21
$ws.Range("T18").Value = 3.9
$ws.Range("W18").Value = 32
$ws.Range("X18").Value = 65
$ws.Range("Z18").Value = 175
$ws.Range("AA18").Value = 65
$ws.Range("AD18").Value = 11.25
$ws.Range("AG18").Value = 11.5
$ws.Range("AP18").Value = 28
$ws.Range("AR18").Value = 150
$ws.Range("AT18").Value = 3.9
$ws.Range("AZ18").Value = 14
$ws.Range("BA18").Value = 28
$ws.Range("BB18").Value = 100
